$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.907.56"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "'1.878.08"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +1.53%  "
$ws.Range("D5").Value = "'334.46"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'1.017"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("D7").Value = "'0.4687"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "'0.3908"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").Value = "'46.84"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").Value = "'0.07942"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").Value = "'1.006"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "'21.58"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").Value = "'1.887.51"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "'5.950"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'7.090"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "'1.020"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "'0.06779"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "'87.55"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'0.00001040"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "'17.02"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").Value = "'27.906.98"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'5.470"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "'2.356"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("D26").Value = "'2.103.73"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").Value = "'159.61"
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("D28").Value = "'19.88"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").Value = "'2.072"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "'5.456"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "'120.86"
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").Value = "'0.09521"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").Value = "'0.9554"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").Value = "'3.657"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").Value = "'5.306"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "'1.351"
$ws.Range("E36").Value = "  -7.07%  "
$ws.Range("D37").Value = "'0.06104"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'0.02237"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").Value = "'1.201"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").Value = "'1.016"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").Value = "'8.103"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").Value = "'0.5883"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("D43").Value = "'0.1889"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "'10.19"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "'1.273"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "'0.5640"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "'12.12"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "'3.385"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "'1.916"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").Value = "'0.06858"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").Value = "'113.50"
$ws.Range("E51").Value = "  +0.93%  "

# The apostrophe above forces Excel to treat the numeric-looking price strings
# as text (matching the source data, which stores prices as plain text - e.g.
# to preserve trailing zeros / thousand-dot formatting). Clear the resulting
# "quote prefix" number format afterwards so cell styling is left untouched.
$ws.Range("D2:D51").ClearFormats()
